# Natmi following Dr Hou advice
# Re-run of the LR-pair analysis for Vwf-Itga2b now includes a third
# sending/target cluster ("FAPs") alongside the existing "ECs" and "sCs",
# giving the full 3x3 cartesian product of clusters (rows 2-10) with
# refreshed expression / specificity statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = New-Object System.Collections.ArrayList
$rows.Add(@("ECs","Vwf","Itga2b","ECs",3,1,58.356725,175.070175,0.9850455027887332,0.9850455027887334,3,1,1.287649333333333,3.862948,0.1870471291026542,0.1870471291026542,75.14299804176667,676.2869823759,0.1842499333321131,0.1842499333321131)) | Out-Null
$rows.Add(@("ECs","Vwf","Itga2b","FAPs",3,1,58.356725,175.070175,0.9850455027887332,0.9850455027887334,3,1,4.622908,13.868724,0.6715350578151914,0.6715350578151914,269.7777708563,2427.9999377067,0.6614925886658262,0.6614925886658264)) | Out-Null
$rows.Add(@("ECs","Vwf","Itga2b","sCs",3,1,58.356725,175.070175,0.9850455027887332,0.9850455027887334,3,1,0.9735329999999999,2.920599,0.1414178130821545,0.1414178130821545,56.81219755942499,511.3097780348249,0.139302980790794,0.139302980790794)) | Out-Null
$rows.Add(@("FAPs","Vwf","Itga2b","ECs",1,0.3333333333333333,0.3120586666666667,0.936176,0.005267464653066949,0.00526746465306695,3,1,1.287649333333333,3.862948,0.1870471291026542,0.1870471291026542,0.4018221340942222,3.616399206848,0.000985264141005881,0.0009852641410058813)) | Out-Null
$rows.Add(@("FAPs","Vwf","Itga2b","FAPs",1,0.3333333333333333,0.3120586666666667,0.936176,0.005267464653066949,0.00526746465306695,3,1,4.622908,13.868724,0.6715350578151914,0.6715350578151914,1.442618506602666,12.983566559424,0.003537287180336791,0.003537287180336791)) | Out-Null
$rows.Add(@("FAPs","Vwf","Itga2b","sCs",1,0.3333333333333333,0.3120586666666667,0.936176,0.005267464653066949,0.00526746465306695,3,1,0.9735329999999999,2.920599,0.1414178130821545,0.1414178130821545,0.303799409936,2.734194689424,0.0007449133317242775,0.0007449133317242776)) | Out-Null
$rows.Add(@("sCs","Vwf","Itga2b","ECs",3,1,0.5738856666666666,1.721657,0.00968703255819983,0.009687032558199832,3,1,1.287649333333333,3.862948,0.1870471291026542,0.1870471291026542,0.7389634960928889,6.650671464836001,0.001811931629535218,0.001811931629535218)) | Out-Null
$rows.Add(@("sCs","Vwf","Itga2b","FAPs",3,1,0.5738856666666666,1.721657,0.00968703255819983,0.009687032558199832,3,1,4.622908,13.868724,0.6715350578151914,0.6715350578151914,2.653020639518666,23.877185755668,0.006505181969028364,0.006505181969028366)) | Out-Null
$rows.Add(@("sCs","Vwf","Itga2b","sCs",3,1,0.5738856666666666,1.721657,0.00968703255819983,0.009687032558199832,3,1,0.9735329999999999,2.920599,0.1414178130821545,0.1414178130821545,0.5586966347269999,5.028269712543,0.001369918959636248,0.001369918959636249)) | Out-Null

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}